# Improve keywords: add a new course entry ("Programming C++" / "C++程式設計")
# at row 15 of the All_EE_Courses sheet, pushing the existing rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("All_EE_Courses")

# Insert a new blank row above the current row 15 (shifts rows 15.. down by one,
# carrying their formatting/styles with them - matches Excel's native "Insert" UI).
$ws.Rows.Item(15).Insert()

# Populate the newly inserted row with the new course (Chinese name in column A,
# English name in column B - same layout as every other row in this sheet).
# (English value is entered first so the shared-string table ends up ordered
# exactly as in the authored workbook.)
$ws.Range("B15").Value = "Programming C++"
$ws.Range("A15").Value = "C++程式設計"

# Leave the view focused on the freshly-added row, same as the author did.
$ws.Range("A15").Select()
